$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2186.45
$ws.Range("I33").Value = 345.23077
$ws.Range("J33").Value = 5605.857
$ws.Range("K33").Value = 345.23077
$ws.Range("L33").Value = 5605.857
$ws.Range("M33").Value = -116.23077
$ws.Range("N33").Value = -6063.857
$ws.Range("H113").Value = 4205.1577
$ws.Range("I113").Value = 3469.8462
$ws.Range("J113").Value = 5798.3335
$ws.Range("K113").Value = 3469.8462
$ws.Range("L113").Value = 5798.3335
$ws.Range("M113").Value = -215.8462
$ws.Range("N113").Value = -12306.3335
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 162.11111
$ws.Range("I5").Value = 179.3
$ws.Range("K5").Value = 179.3
$ws.Range("M5").Value = -67.30000000000001
$ws.Range("H32").Value = 4296.0107
$ws.Range("I32").Value = 1988.0588
$ws.Range("J32").Value = 23913.6
$ws.Range("K32").Value = 1988.0588
$ws.Range("L32").Value = 23913.6
$ws.Range("M32").Value = -1701.0588
$ws.Range("N32").Value = -24487.6
$ws.Range("H35").Value = 4493.5
$ws.Range("I35").Value = 4493.5
$ws.Range("K35").Value = 4493.5
$ws.Range("M35").Value = -4087.5
$ws.Range("H39").Value = 57841.668
$ws.Range("J39").Value = 81762.5
$ws.Range("L39").Value = 81762.5
$ws.Range("N39").Value = -82802.5
$ws.Range("H61").Value = 8168.6875
$ws.Range("I61").Value = 8799.909
$ws.Range("J61").Value = 6780
$ws.Range("K61").Value = 8799.909
$ws.Range("L61").Value = 6780
$ws.Range("M61").Value = -8587.909
$ws.Range("N61").Value = -7204
$ws.Range("H74").Value = 3343.575
$ws.Range("I74").Value = 1372.5358
$ws.Range("J74").Value = 7942.6665
$ws.Range("K74").Value = 1372.5358
$ws.Range("L74").Value = 7942.6665
$ws.Range("M74").Value = -498.5358000000001
$ws.Range("N74").Value = -9690.666499999999
$ws.Range("H77").Value = 3343.575
$ws.Range("I77").Value = 1372.5358
$ws.Range("J77").Value = 7942.6665
$ws.Range("K77").Value = 6862.679
$ws.Range("L77").Value = 39713.3325
$ws.Range("M77").Value = -2494.679
$ws.Range("N77").Value = -48449.3325
$ws.Range("H88").Value = 2209
$ws.Range("I88").Value = 1646.8
$ws.Range("J88").Value = 2677.5
$ws.Range("K88").Value = 1646.8
$ws.Range("L88").Value = 2677.5
$ws.Range("M88").Value = -1240.8
$ws.Range("N88").Value = -3489.5
$ws.Range("H91").Value = 2209
$ws.Range("I91").Value = 1646.8
$ws.Range("J91").Value = 2677.5
$ws.Range("K91").Value = 1646.8
$ws.Range("L91").Value = 2677.5
$ws.Range("M91").Value = -242.8
$ws.Range("N91").Value = -5485.5
$ws.Range("H102").Value = 3742.8
$ws.Range("I102").Value = 3772.0715
$ws.Range("K102").Value = 3772.0715
$ws.Range("M102").Value = -2150.0715
$ws.Range("H122").Value = 3178.2903
$ws.Range("I122").Value = 2039.5385
$ws.Range("K122").Value = 6118.6155
$ws.Range("M122").Value = -3668.6155
$ws.Range("H132").Value = 4649.3706
$ws.Range("I132").Value = 2291.4707
$ws.Range("J132").Value = 8657.799999999999
$ws.Range("K132").Value = 6874.4121
$ws.Range("L132").Value = 25973.4
$ws.Range("M132").Value = -4344.4121
$ws.Range("N132").Value = -31033.4
$ws.Range("H136").Value = 8168.6875
$ws.Range("I136").Value = 8799.909
$ws.Range("J136").Value = 6780
$ws.Range("K136").Value = 26399.727
$ws.Range("L136").Value = 20340
$ws.Range("M136").Value = -23849.727
$ws.Range("N136").Value = -25440
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 162.11111
$ws.Range("I4").Value = 179.3
$ws.Range("K4").Value = 179.3
$ws.Range("M4").Value = -64.30000000000001
$ws.Range("H54").Value = 60749.25
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 60749.25
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 60749.25
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -61717.25
$ws.Range("H94").Value = 1055.6923
$ws.Range("I94").Value = 781.7727
$ws.Range("K94").Value = 781.7727
$ws.Range("M94").Value = -330.7727
$ws.Range("H99").Value = 3107.0715
$ws.Range("I99").Value = 2999.923
$ws.Range("K99").Value = 2999.923
$ws.Range("M99").Value = -1501.923
$ws.Range("H105").Value = 3317.5789
$ws.Range("I105").Value = 3277.3572
$ws.Range("K105").Value = 3277.3572
$ws.Range("M105").Value = -1530.3572
$ws.Range("H107").Value = 2859.5483
$ws.Range("I107").Value = 2945.8
$ws.Range("K107").Value = 2945.8
$ws.Range("M107").Value = -1025.8
$ws.Range("H134").Value = 5352.857
$ws.Range("I134").Value = 2967.2727
$ws.Range("K134").Value = 8901.8181
$ws.Range("M134").Value = -6366.8181
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1164.3125
$ws.Range("I16").Value = 862.1667
$ws.Range("J16").Value = 2070.75
$ws.Range("K16").Value = 862.1667
$ws.Range("L16").Value = 2070.75
$ws.Range("M16").Value = -575.1667
$ws.Range("N16").Value = -2644.75
$ws.Range("H86").Value = 5912.273
$ws.Range("I86").Value = 6073.5
$ws.Range("J86").Value = 5718.8
$ws.Range("K86").Value = 6073.5
$ws.Range("L86").Value = 5718.8
$ws.Range("M86").Value = -4950.5
$ws.Range("N86").Value = -7964.8
$ws.Range("H89").Value = 5912.273
$ws.Range("I89").Value = 6073.5
$ws.Range("J89").Value = 5718.8
$ws.Range("K89").Value = 30367.5
$ws.Range("L89").Value = 28594
$ws.Range("M89").Value = -24751.5
$ws.Range("N89").Value = -39826
$ws.Range("H105").Value = 3580.7
$ws.Range("I105").Value = 3213.375
$ws.Range("K105").Value = 3213.375
$ws.Range("M105").Value = -1466.375
$ws.Range("H107").Value = 1403.1052
$ws.Range("I107").Value = 1229.7142
$ws.Range("J107").Value = 1888.6
$ws.Range("K107").Value = 1229.7142
$ws.Range("L107").Value = 1888.6
$ws.Range("M107").Value = 690.2858000000001
$ws.Range("N107").Value = -5728.6
$ws.Range("H113").Value = 1164.3125
$ws.Range("I113").Value = 862.1667
$ws.Range("J113").Value = 2070.75
$ws.Range("K113").Value = 862.1667
$ws.Range("L113").Value = 2070.75
$ws.Range("M113").Value = 1307.8333
$ws.Range("N113").Value = -6410.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1082.25
$ws.Range("I2").Value = 176.25
$ws.Range("K2").Value = 1057.5
$ws.Range("M2").Value = -944.5
$ws.Range("H17").Value = 594
$ws.Range("I17").Value = 546.3333
$ws.Range("J17").Value = 880
$ws.Range("K17").Value = 1638.9999
$ws.Range("L17").Value = 2640
$ws.Range("M17").Value = -1469.9999
$ws.Range("N17").Value = -2978
$ws.Range("H34").Value = 2160.8333
$ws.Range("J34").Value = 3526
$ws.Range("L34").Value = 10578
$ws.Range("N34").Value = -10746
$ws.Range("H39").Value = 4370.75
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H55").Value = 3136.75
$ws.Range("I55").Value = 865.6667
$ws.Range("J55").Value = 9950
$ws.Range("K55").Value = 2597.0001
$ws.Range("L55").Value = 29850
$ws.Range("M55").Value = -2420.0001
$ws.Range("N55").Value = -30204
$ws.Range("H128").Value = 449998.75
$ws.Range("I128").Value = 449998.75
$ws.Range("K128").Value = 1349996.25
$ws.Range("M128").Value = -1345016.25
$ws.Range("H137").Value = 3307.8823
$ws.Range("I137").Value = 2110.8
$ws.Range("J137").Value = 3806.6667
$ws.Range("K137").Value = 6332.400000000001
$ws.Range("L137").Value = 11420.0001
$ws.Range("M137").Value = -1232.400000000001
$ws.Range("N137").Value = -21620.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6478.8423
$ws.Range("I70").Value = 5875
$ws.Range("J70").Value = 6549.8823
$ws.Range("K70").Value = 5875
$ws.Range("L70").Value = 6549.8823
$ws.Range("M70").Value = -5605
$ws.Range("N70").Value = -7089.8823
$ws.Range("H73").Value = 6478.8423
$ws.Range("I73").Value = 5875
$ws.Range("J73").Value = 6549.8823
$ws.Range("K73").Value = 5875
$ws.Range("L73").Value = 6549.8823
$ws.Range("M73").Value = -4939
$ws.Range("N73").Value = -8421.882300000001
$ws.Range("H113").Value = 4314.773
$ws.Range("I113").Value = 3531.8572
$ws.Range("J113").Value = 5684.875
$ws.Range("K113").Value = 3531.8572
$ws.Range("L113").Value = 5684.875
$ws.Range("M113").Value = -1361.8572
$ws.Range("N113").Value = -10024.875
$ws.Range("H126").Value = 4866.467
$ws.Range("I126").Value = 2999.7273
$ws.Range("K126").Value = 8999.1819
$ws.Range("M126").Value = -6529.1819
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 94914.766
$ws.Range("I7").Value = 141737.5
$ws.Range("K7").Value = 141737.5
$ws.Range("M7").Value = -141625.5
$ws.Range("H55").Value = 782.55554
$ws.Range("I55").Value = 383.30768
$ws.Range("J55").Value = 1153.2858
$ws.Range("K55").Value = 383.30768
$ws.Range("L55").Value = 1153.2858
$ws.Range("M55").Value = -210.30768
$ws.Range("N55").Value = -1499.2858
$ws.Range("H122").Value = 6280.826
$ws.Range("I122").Value = 5313.684
$ws.Range("J122").Value = 10874.75
$ws.Range("K122").Value = 15941.052
$ws.Range("L122").Value = 32624.25
$ws.Range("M122").Value = -13491.052
$ws.Range("N122").Value = -37524.25
$ws.Range("H126").Value = 94914.766
$ws.Range("I126").Value = 141737.5
$ws.Range("K126").Value = 425212.5
$ws.Range("M126").Value = -422742.5
$ws.Range("H132").Value = 4054.05
$ws.Range("I132").Value = 2395.2144
$ws.Range("J132").Value = 7924.6665
$ws.Range("K132").Value = 7185.6432
$ws.Range("L132").Value = 23773.9995
$ws.Range("M132").Value = -4655.6432
$ws.Range("N132").Value = -28833.9995
$ws.Range("H136").Value = 4529.222
$ws.Range("I136").Value = 2443.1428
$ws.Range("J136").Value = 5856.727
$ws.Range("K136").Value = 7329.428400000001
$ws.Range("L136").Value = 17570.181
$ws.Range("M136").Value = -4779.428400000001
$ws.Range("N136").Value = -22670.181
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 18000
$ws.Range("J58").Value = 18000
$ws.Range("L58").Value = 18000
$ws.Range("N58").Value = -18616
$ws.Range("H136").Value = 13034
$ws.Range("I136").Value = 14735
$ws.Range("J136").Value = 11333
$ws.Range("K136").Value = 44205
$ws.Range("L136").Value = 33999
$ws.Range("M136").Value = -41655
$ws.Range("N136").Value = -39099

Write-Output "applied 269 cell updates"